$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 1193
$wsExhibit.Range("F6").Value = 72
$wsExhibit.Range("F7").Value = 4375
$wsExhibit.Range("F8").Value = 2604
$wsExhibit.Range("F10").Value = 2520
$wsExhibit.Range("F14").Value = 1663
$wsExhibit.Range("F15").Value = 664
$wsExhibit.Range("F17").Value = 113
$wsExhibit.Range("F18").Value = 330
$wsExhibit.Range("F19").Value = 27
$wsExhibit.Range("F22").Value = 31
$wsExhibit.Range("F23").Value = 482
$wsExhibit.Range("F26").Value = 553
$wsExhibit.Range("F27").Value = 694
$wsExhibit.Range("F30").Value = 409
$wsExhibit.Range("F32").Value = 1619
$wsExhibit.Range("F33").Value = 1040
$wsExhibit.Range("F34").Value = 135
$wsExhibit.Range("F36").Value = 1141
$wsExhibit.Range("F37").Value = 2057
$wsExhibit.Range("F38").Value = 271
$wsExhibit.Range("F40").Value = 547
$wsExhibit.Range("F43").Value = 660
$wsExhibit.Range("F44").Value = 1330
$wsExhibit.Range("F45").Value = 105
$wsExhibit.Range("F47").Value = 436
$wsExhibit.Range("F48").Value = 70

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 11

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1193
$wsAll.Range("F4").Value = 11
$wsAll.Range("F5").Value = 72
$wsAll.Range("F6").Value = 4375
$wsAll.Range("F7").Value = 2604
$wsAll.Range("F8").Value = 2520
$wsAll.Range("F9").Value = 1663
$wsAll.Range("F12").Value = 664
$wsAll.Range("F14").Value = 113
$wsAll.Range("F15").Value = 330
$wsAll.Range("F16").Value = 27
$wsAll.Range("F19").Value = 482
$wsAll.Range("F22").Value = 553
$wsAll.Range("F23").Value = 694
$wsAll.Range("F29").Value = 409
$wsAll.Range("F30").Value = 1619
$wsAll.Range("F31").Value = 1040
$wsAll.Range("F32").Value = 135
$wsAll.Range("F35").Value = 2057
$wsAll.Range("F36").Value = 271
$wsAll.Range("F40").Value = 547
$wsAll.Range("F43").Value = 660
$wsAll.Range("F44").Value = 1330
$wsAll.Range("F46").Value = 105
$wsAll.Range("F47").Value = 436
$wsAll.Range("F48").Value = 70
